# Update the Gantt chart "Installation" group (rows 52-57) with the new
# "Item Manager" task breakdown, per the commit:
# "Update Gantt Chart with Item Manager Tasks"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# Replace the placeholder "TBD" rows (A52:A57) with the real task names,
# and fill in hours spent on each (column B).
$ws.Range("A52").Value = "Requirements Collection"
$ws.Range("B52").Value = 3

$ws.Range("A53").Value = "Function Definitions"
$ws.Range("B53").Value = 2

$ws.Range("A54").Value = "User Documentation"
$ws.Range("B54").Value = 1

$ws.Range("A55").Value = "Programming"
$ws.Range("B55").Value = 3

$ws.Range("A56").Value = "Testing"
$ws.Range("B56").Value = 2

$ws.Range("A57").Value = "Installation"
$ws.Range("B57").Value = 2

$excel.CalculateFullRebuild()
